$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (header + 16 player rows) in the new order, with the
# "Kris Dunn" row removed entirely.
$data = @(
    @("Oyuncu Adı", "Pozisyon", "Takım"),
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("T.J. McConnell", "PG", "Indiana Pacers"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Joel Embiid", "C", "Philadelphia 76ers")
)

$oldLastRow = 18
$newLastRow = $data.Count

# Remove the row that is no longer needed (the table shrinks by one row)
# so everything shifts up cleanly and no stale cells are left behind.
if ($newLastRow -lt $oldLastRow) {
    $rowsToDelete = $oldLastRow - $newLastRow
    for ($i = 0; $i -lt $rowsToDelete; $i++) {
        $ws.Rows.Item($oldLastRow).Delete()
    }
}

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}
